$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, pushing the old rows 4 and 5 down to 5 and 6.
$ws.Rows.Item(4).Insert()

# Copy style of date cell (D) from the row below (old row 4, now row 5) to the new D4.
$ws.Range("D5").Copy()
$ws.Range("D4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new weekly row of data.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44784
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112012
$ws.Range("G4").Value = "Espinaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("N4").Value = "`$/cuna 10 kilos"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 850
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"
